$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("IVV", 432.7099914550781, 66.1177978515625, 21.71),
    @("VOO", 395.7799987792969, 60.5667724609375, 21.56),
    @("VTI", 213.0099945068359, 29.48443603515625, 20.03),
    @("QQQ", 363.4400024414062, 104.7037963867188, 29.64),
    @("VTV", 137.3999938964844, 5.09356689453125, 15.36),
    @("VUG", 276.8900146484375, 69.81861877441406, 33.33),
    @("IWF", 271.0299987792969, 62.83857727050781, 29.92),
    @("VIG", 156.5899963378906, 15.28024291992188, 22.04),
    @("IJH", 242.5599975585938, 10.45112609863281, 13.36),
    @("IJR", 91.66000366210938, -1.283927917480469, 11.33)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $row++
}
